$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "[PyTorch] 모델의 구조도 요약(summary) 출력 (torchsummary)"
$ws.Range("E4").Value = "https://teddylee777.github.io/pytorch/pytorch-torchsummary"

$ws.Range("D9").Value = "[공지] SIAI 커리큘럼 관련 소개"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/notice-curriculum-details/#utm_source=rss&utm_medium=rss&utm_campaign=notice-curriculum-details"

$ws.Range("D20").Value = "[책] [AI/MLOps] 쿠브플루 운영 가이드 (Kubeflow Operations Guide)"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/610"

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D42").Value = "Boost 라이브러리 환경설정 및 Visual Studio 2019 디렉토리 설정"
$ws.Range("E42").Value = "https://kjk92.tistory.com/81"

$ws.Range("D50").Value = "be a catch-22 situation"
$ws.Range("E50").Value = "http://incredible.egloos.com/7535716"
